$d = $word.ActiveDocument

# --- 1) Split the "Secondary LIA" run into two runs at the same boundary ---
# "Secondary LIA – mixed with CT2 signal"  ->
# "Secondary LIA – mixed " | "with CT2 signal"
$r = $d.Content
$found = $r.Find.Execute("with CT2 signal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'with CT2 signal' text"
}
# Toggling a character property on just this sub-range forces Word to
# materialize it as its own run (the earlier text keeps its own run too).
$r.Bold = 1
$r.Bold = 0

# --- 2) Add list numbering (numId 1, ilvl 0) to the "Bx/By reference" paragraph ---
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Bx/By reference*output of the*Mux*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the 'Bx/By reference' paragraph"
}

# Reuse the same list (numId 1) already used by the preceding "Secondary LIA"
# bullet, rather than minting a brand-new list/abstractNum.
$sourceList = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Secondary LIA*mixed*CT2 signal*") {
        $sourceList = $p
        break
    }
}
if ($sourceList -eq $null) {
    throw "Could not find the 'Secondary LIA' paragraph"
}

$tmpl = $sourceList.Range.ListFormat.ListTemplate
$target.Range.ListFormat.ApplyListTemplateWithLevel($tmpl, $true)
